$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 303, pushing existing rows 303-319 down to 304-320
$ws.Rows(303).Insert()

# Populate the newly inserted row 303 with the new record's data
$ws.Range("A303").Value = 7
$ws.Range("B303").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C303").Value = "Ñuble"
$ws.Range("D303").Value = 44585
$ws.Range("E303").Value = 16
$ws.Range("F303").Value = 100114014
$ws.Range("G303").Value = "Betarraga"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 400
$ws.Range("K303").Value = 600
$ws.Range("L303").Value = 650
$ws.Range("M303").Value = 625
$ws.Range("N303").Value = "$/paquete 5 unidades"
$ws.Range("O303").Value = "Región del Maule"
$ws.Range("P303").Value = 125
$ws.Range("Q303").Value = 5
$ws.Range("R303").Value = "Hortaliza"
